# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Column D ("Price") cells are stored as plain text in this workbook (values
# like "3.151.36" use '.' as a thousands separator and aren't valid numbers).
# Writing such text through .Value lets Excel's input parser reinterpret it
# as a real number (dropping significant trailing zeros, mis-parsing the
# grouping dots, etc.), so those writes are wrapped with a temporary "@"
# (text) number format and reset back to the default "Normal" style
# afterwards so no stray explicit format is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '57.844.00'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.13%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.152.40'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.40%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '531.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.45%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '140.71'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.17%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.151.78'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.39%  '

# Row 9
$ws.Range('E9').Value = '  +2.77%  '

# Row 10
$ws.Range('E10').Value = '  -0.58%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.109'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.15%  '

# Row 12
$ws.Range('E12').Value = '  +4.53%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.693.24'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.37%  '

# Row 14
$ws.Range('E14').Value = '  +2.85%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '25.64'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.15%  '

# Row 16
$ws.Range('E16').Value = '  +0.63%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '58.004.08'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.26%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.145.02'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.18%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.13'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.44%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.84'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.17%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.00'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.77%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '355.17'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +5.17%  '

# Row 23
$ws.Range('E23').Value = '  +0.05%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.61'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.09%  '

# Row 25
$ws.Range('E25').Value = '  +0.60%  '

# Row 26
$ws.Range('E26').Value = '  +1.47%  '

# Row 27
$ws.Range('E27').Value = '  -0.04%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0941'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.01%  '

# Row 29
$ws.Range('E29').Value = '  +3.20%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.08%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.39'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.57%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.90'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.16%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '21.29'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.60%  '

# Row 34
$ws.Range('E34').Value = '  +0.58%  '

# Row 35
$ws.Range('E35').Value = '  +5.68%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '157.66'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.23%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.20'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.89%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '26.18'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.28%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.28'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.04%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0672'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.74%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.62'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +10.64%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.09'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.05%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.705'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.00%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.193.14'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.27%  '

# Row 45
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0273'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +5.79%  '

# Row 46
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '36.72'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.64%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.13%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.335.48'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.41%  '

# Row 49
$ws.Range('E49').Value = '  +2.41%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.07'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.20%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '20.39'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.36%  '
